# Update "想去人数" (want-to-go headcount) figures that were refreshed
# when the gh-pages data was regenerated (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# -------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 26940
$ws1.Range("F4").Value  = 600
$ws1.Range("F7").Value  = 181
$ws1.Range("F10").Value = 368
$ws1.Range("F11").Value = 460
$ws1.Range("F12").Value = 193
$ws1.Range("F13").Value = 52
$ws1.Range("F14").Value = 306
$ws1.Range("F15").Value = 86
$ws1.Range("F16").Value = 452
$ws1.Range("F17").Value = 63
$ws1.Range("F19").Value = 228
$ws1.Range("F20").Value = 65

# -------------------------------------------------------------------
# Sheet "演出" (Performances)
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 114

# -------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# -------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5129
$ws3.Range("F3").Value = 255

# -------------------------------------------------------------------
# Sheet "全部类型" (All types - aggregate of the sheets above)
# -------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5129
$ws4.Range("F4").Value  = 255
$ws4.Range("F5").Value  = 26940
$ws4.Range("F6").Value  = 600
$ws4.Range("F13").Value = 181
$ws4.Range("F17").Value = 114
$ws4.Range("F22").Value = 368
$ws4.Range("F23").Value = 460
$ws4.Range("F24").Value = 193
$ws4.Range("F25").Value = 52
$ws4.Range("F27").Value = 306
$ws4.Range("F28").Value = 86
$ws4.Range("F31").Value = 452
$ws4.Range("F32").Value = 63
$ws4.Range("F35").Value = 228
$ws4.Range("F37").Value = 65
